$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.416099999999992
$ws.Range("B4").Value = 4.795600000000001
$ws.Range("C6").Value = -10.8008
$ws.Range("B7").Value = 5.813900000000001
$ws.Range("C7").Value = -11.2971
$ws.Range("B8").Value = 5.110399999999997
$ws.Range("C8").Value = -10.58839999999999
$ws.Range("A11").Value = -22.03900000000002
$ws.Range("D11").Value = -8.890199999999988
$ws.Range("A12").Value = -22.72950000000001
$ws.Range("B12").Value = 5.2123
$ws.Range("B14").Value = 9.733100000000006
$ws.Range("D14").Value = -8.811499999999997
$ws.Range("A15").Value = -21.43430000000002
$ws.Range("C19").Value = -12.6415
$ws.Range("D19").Value = -8.638199999999987
$ws.Range("C21").Value = -12.5359
$ws.Range("D21").Value = -8.76949999999999
$ws.Range("B22").Value = 4.968300000000004
$ws.Range("C24").Value = -11.5334
$ws.Range("C25").Value = -10.61169999999999
